# Update header labels on row 1 so the first row can be used automatically
# as a header when the data is loaded into Power BI.

$wb = $excel.ActiveWorkbook

# Sheets that use "Ano" (Year) prefixed labels for columns B1:E1
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet that uses "Intervalo" prefixed labels for columns B1:E1
$ws = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws.Range("B1").Value = "Intervalo 2015"
$ws.Range("C1").Value = "Intervalo 2015-2030"
$ws.Range("D1").Value = "Intervalo 2031-2040"
$ws.Range("E1").Value = "Intervalo 2041-2050"

# Sheet with only a single year column (B1)
$ws = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws.Range("B1").Value = "Ano 2015"
